$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the original (pre-edit) values for rows 21-23 before any writes ---
# These three rows cyclically rotate their "species record" data:
#   new row21 <- old row22
#   new row22 <- old row23
#   new row23 <- old row21

$row21_A = $ws.Range("A21").Value()
$row21_B = $ws.Range("B21").Value()
$row21_E = $ws.Range("E21").Value()
$row21_F = $ws.Range("F21").Value()
$row21_G = $ws.Range("G21").Value()
$row21_H = $ws.Range("H21").Value()

$row22_A = $ws.Range("A22").Value()
$row22_B = $ws.Range("B22").Value()
$row22_E = $ws.Range("E22").Value()
$row22_F = $ws.Range("F22").Value()
$row22_G = $ws.Range("G22").Value()
$row22_H = $ws.Range("H22").Value()

$row23_A = $ws.Range("A23").Value()
$row23_B = $ws.Range("B23").Value()
$row23_E = $ws.Range("E23").Value()
$row23_F = $ws.Range("F23").Value()
$row23_G = $ws.Range("G23").Value()
$row23_H = $ws.Range("H23").Value()

# =========================================================================
# ROW 21: becomes old row22's species record; Q/R rounded to whole metres;
# Starttid/Sluttid ("00:00") cells are cleared (removed).
# =========================================================================
$ws.Range("A21").Value = $row22_A
$ws.Range("B21").Value = $row22_B
$ws.Range("E21").Value = $row22_E
$ws.Range("F21").Value = $row22_F
$ws.Range("G21").Value = $row22_G
$ws.Range("H21").Value = $row22_H

$ws.Range("Q21").Value = 383386
$ws.Range("R21").Value = 6664494

$ws.Range("Z21").ClearContents()
$ws.Range("AB21").ClearContents()

# =========================================================================
# ROW 22: becomes old row23's species record (a bird observation), so the
# Antal/Enhet/Aktivitet columns take the bird-record shape: I22 gets a
# count, J22 (Enhet) is cleared, L22 (Kon) becomes an empty-but-present
# cell, M22 (Aktivitet) gets text, and AF22 (Bestamningsmetod) is cleared.
# =========================================================================
$ws.Range("A22").Value = $row23_A
$ws.Range("B22").Value = $row23_B
$ws.Range("E22").Value = $row23_E
$ws.Range("F22").Value = $row23_F
$ws.Range("G22").Value = $row23_G
$ws.Range("H22").Value = $row23_H

$ws.Range("I22").NumberFormat = "@"
$ws.Range("I22").Value = "1"
$ws.Range("I22").Style = "Normal"
$ws.Range("J22").ClearContents()
$ws.Range("L22").Style = "Normal"
$ws.Range("M22").Value = "lockläte, övriga läten"

$ws.Range("Q22").Value = 383215
$ws.Range("R22").Value = 6664539
$ws.Range("S22").Value = 25

$ws.Range("Z22").ClearContents()
$ws.Range("AB22").ClearContents()
$ws.Range("AF22").ClearContents()

# =========================================================================
# ROW 23: becomes old row21's species record (a fungus observation), so the
# Antal column is cleared back out, J23 (Enhet) becomes an empty-but-
# present cell again, L23/M23 (Kon/Aktivitet) are cleared, and AF23
# (Bestamningsmetod) becomes an empty-but-present cell.
# =========================================================================
$ws.Range("A23").Value = $row21_A
$ws.Range("B23").Value = $row21_B
$ws.Range("E23").Value = $row21_E
$ws.Range("F23").Value = $row21_F
$ws.Range("G23").Value = $row21_G
$ws.Range("H23").Value = $row21_H

$ws.Range("I23").ClearContents()
$ws.Range("J23").Style = "Normal"
$ws.Range("L23").ClearContents()
$ws.Range("M23").ClearContents()

$ws.Range("Q23").Value = 383311
$ws.Range("R23").Value = 6664460
$ws.Range("S23").Value = 10

$ws.Range("Z23").ClearContents()
$ws.Range("AB23").ClearContents()
$ws.Range("AF23").Style = "Normal"

# =========================================================================
# ROW 24: same observation, only the coordinates are rounded to whole
# metres and the Starttid/Sluttid cells are cleared.
# =========================================================================
$ws.Range("Q24").Value = 383318
$ws.Range("R24").Value = 6664423

$ws.Range("Z24").ClearContents()
$ws.Range("AB24").ClearContents()
